$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1748463629766/S84YrFTsc0-camera-video-segment-driver-1748463632266.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOFTLUGLO%2F20250529%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250529T210320Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEMv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIB2ZTYdObBkz97dPlxhBdWMP%2B%2FWDa8OFSZfyYXQA2ZzgAiA4g56GZxpBRUH4Q4nFWXa6dLWjwViVIhrBPNF6wT0kwirmAwiU%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMuezOu6xnbDig%2BzMOKroD7cPvdDW9pey5Wwt65aZeLUFzA%2FU17M23Vg4GRbV8IjitQm87hihB6MDN%2Fy6kf2QIngVv%2FrPBjQVDos0mRRLFr0oxpIm%2BcLBJKEuUqC9oE%2BnO5oWndifBX1XYlPv688N5yPfu2hpk9D7zxv%2FyEoZLzQP6l3woQXUzYn7xDSlNsHy36DzswM%2FmHAy6tIDwMaI9fz%2BNUs0NdnoYceo8AQIkKwoy4qRmX%2FDVwhVEmazuPIUn4ZP096AlnUDQQr41zrh1bKzzj4qckTCUZTbIOT9VoF%2B02gy5IhJPG8kSvI3Ybrfnz0N%2Ba5Dq6t7INTl1BPr3Ja6Oto31MzTIh9zh6DZvwFIyHbIbSh44iK3zjhYG7JoBQZVHXHdJVEs%2FQmBnONqJoAmGWT8DinWK7v26dUe%2BSUOxDQGRJrGIcOoDlaAlPDHZSF3c2MCtfk4guCBVIdr%2BuJrzSqdG2zwbo4TSopFoXhcLDzf2zNfWJygj9eyvgp9%2FBpBkZhA7gWJdvMXebxN8h9YSUVZ6vjXo5ZP7MzL6YlaIJj6txOnMvOVcxcnZLwayndknFByWtxfp8p42nBOSCJnnTIfFhhVRhDCI3uLBBjqmAaThJMX71dNkCv3kni%2B3hFSEJnFJdS1RyXwWlaHQtiEKWo8ugd0oAhTN1Qij7mXIcqbaoft%2BqxQYS33yozsS1FqNi2pOd5LQWbg3I83LJlgz5J0cmFf9NDMUdCC2Q5DoWGCrkgbni9doXX3XCfPdtcfsbtvFI9F4OJcRk8JeFlFA%2Bagtk54Tb9MIf0UlFiqRjPt0kNARFwLFZ%2FJL366fRNrIjh%2Bs%2B0E%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2030%20May%202025%2005%3A03%3A20%20GMT&X-Amz-Signature=53b05ce4a83307816b167d61d3d8e2d47c7b07be5fe8d0c1df87bc8b22d7c91b"
$ws.Range("L2").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1748463629766/ElYsNGhClW-camera-video-segment-1748463632266.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOFTLUGLO%2F20250529%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250529T210320Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEMv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIB2ZTYdObBkz97dPlxhBdWMP%2B%2FWDa8OFSZfyYXQA2ZzgAiA4g56GZxpBRUH4Q4nFWXa6dLWjwViVIhrBPNF6wT0kwirmAwiU%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMuezOu6xnbDig%2BzMOKroD7cPvdDW9pey5Wwt65aZeLUFzA%2FU17M23Vg4GRbV8IjitQm87hihB6MDN%2Fy6kf2QIngVv%2FrPBjQVDos0mRRLFr0oxpIm%2BcLBJKEuUqC9oE%2BnO5oWndifBX1XYlPv688N5yPfu2hpk9D7zxv%2FyEoZLzQP6l3woQXUzYn7xDSlNsHy36DzswM%2FmHAy6tIDwMaI9fz%2BNUs0NdnoYceo8AQIkKwoy4qRmX%2FDVwhVEmazuPIUn4ZP096AlnUDQQr41zrh1bKzzj4qckTCUZTbIOT9VoF%2B02gy5IhJPG8kSvI3Ybrfnz0N%2Ba5Dq6t7INTl1BPr3Ja6Oto31MzTIh9zh6DZvwFIyHbIbSh44iK3zjhYG7JoBQZVHXHdJVEs%2FQmBnONqJoAmGWT8DinWK7v26dUe%2BSUOxDQGRJrGIcOoDlaAlPDHZSF3c2MCtfk4guCBVIdr%2BuJrzSqdG2zwbo4TSopFoXhcLDzf2zNfWJygj9eyvgp9%2FBpBkZhA7gWJdvMXebxN8h9YSUVZ6vjXo5ZP7MzL6YlaIJj6txOnMvOVcxcnZLwayndknFByWtxfp8p42nBOSCJnnTIfFhhVRhDCI3uLBBjqmAaThJMX71dNkCv3kni%2B3hFSEJnFJdS1RyXwWlaHQtiEKWo8ugd0oAhTN1Qij7mXIcqbaoft%2BqxQYS33yozsS1FqNi2pOd5LQWbg3I83LJlgz5J0cmFf9NDMUdCC2Q5DoWGCrkgbni9doXX3XCfPdtcfsbtvFI9F4OJcRk8JeFlFA%2Bagtk54Tb9MIf0UlFiqRjPt0kNARFwLFZ%2FJL366fRNrIjh%2Bs%2B0E%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2030%20May%202025%2005%3A03%3A20%20GMT&X-Amz-Signature=3d69669a2a8aaf202367472ef6c9565b4359c3a0110d5c3b6b77753f1d4fea6a"
$ws.Range("K3").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990870452/1748461966091/lLthPW38mJ-camera-video-segment-driver-1748461971091.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOFTLUGLO%2F20250529%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250529T210320Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEMv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIB2ZTYdObBkz97dPlxhBdWMP%2B%2FWDa8OFSZfyYXQA2ZzgAiA4g56GZxpBRUH4Q4nFWXa6dLWjwViVIhrBPNF6wT0kwirmAwiU%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMuezOu6xnbDig%2BzMOKroD7cPvdDW9pey5Wwt65aZeLUFzA%2FU17M23Vg4GRbV8IjitQm87hihB6MDN%2Fy6kf2QIngVv%2FrPBjQVDos0mRRLFr0oxpIm%2BcLBJKEuUqC9oE%2BnO5oWndifBX1XYlPv688N5yPfu2hpk9D7zxv%2FyEoZLzQP6l3woQXUzYn7xDSlNsHy36DzswM%2FmHAy6tIDwMaI9fz%2BNUs0NdnoYceo8AQIkKwoy4qRmX%2FDVwhVEmazuPIUn4ZP096AlnUDQQr41zrh1bKzzj4qckTCUZTbIOT9VoF%2B02gy5IhJPG8kSvI3Ybrfnz0N%2Ba5Dq6t7INTl1BPr3Ja6Oto31MzTIh9zh6DZvwFIyHbIbSh44iK3zjhYG7JoBQZVHXHdJVEs%2FQmBnONqJoAmGWT8DinWK7v26dUe%2BSUOxDQGRJrGIcOoDlaAlPDHZSF3c2MCtfk4guCBVIdr%2BuJrzSqdG2zwbo4TSopFoXhcLDzf2zNfWJygj9eyvgp9%2FBpBkZhA7gWJdvMXebxN8h9YSUVZ6vjXo5ZP7MzL6YlaIJj6txOnMvOVcxcnZLwayndknFByWtxfp8p42nBOSCJnnTIfFhhVRhDCI3uLBBjqmAaThJMX71dNkCv3kni%2B3hFSEJnFJdS1RyXwWlaHQtiEKWo8ugd0oAhTN1Qij7mXIcqbaoft%2BqxQYS33yozsS1FqNi2pOd5LQWbg3I83LJlgz5J0cmFf9NDMUdCC2Q5DoWGCrkgbni9doXX3XCfPdtcfsbtvFI9F4OJcRk8JeFlFA%2Bagtk54Tb9MIf0UlFiqRjPt0kNARFwLFZ%2FJL366fRNrIjh%2Bs%2B0E%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2030%20May%202025%2005%3A03%3A20%20GMT&X-Amz-Signature=59ef67a1187dbd53d1723c9eedccc5d331a68ace0d322f2f8123b08cdaeeaad7"
$ws.Range("L3").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474990870452/1748461966091/M27lrHIJ6l-camera-video-segment-1748461971091.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOFTLUGLO%2F20250529%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250529T210320Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEMv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIB2ZTYdObBkz97dPlxhBdWMP%2B%2FWDa8OFSZfyYXQA2ZzgAiA4g56GZxpBRUH4Q4nFWXa6dLWjwViVIhrBPNF6wT0kwirmAwiU%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMuezOu6xnbDig%2BzMOKroD7cPvdDW9pey5Wwt65aZeLUFzA%2FU17M23Vg4GRbV8IjitQm87hihB6MDN%2Fy6kf2QIngVv%2FrPBjQVDos0mRRLFr0oxpIm%2BcLBJKEuUqC9oE%2BnO5oWndifBX1XYlPv688N5yPfu2hpk9D7zxv%2FyEoZLzQP6l3woQXUzYn7xDSlNsHy36DzswM%2FmHAy6tIDwMaI9fz%2BNUs0NdnoYceo8AQIkKwoy4qRmX%2FDVwhVEmazuPIUn4ZP096AlnUDQQr41zrh1bKzzj4qckTCUZTbIOT9VoF%2B02gy5IhJPG8kSvI3Ybrfnz0N%2Ba5Dq6t7INTl1BPr3Ja6Oto31MzTIh9zh6DZvwFIyHbIbSh44iK3zjhYG7JoBQZVHXHdJVEs%2FQmBnONqJoAmGWT8DinWK7v26dUe%2BSUOxDQGRJrGIcOoDlaAlPDHZSF3c2MCtfk4guCBVIdr%2BuJrzSqdG2zwbo4TSopFoXhcLDzf2zNfWJygj9eyvgp9%2FBpBkZhA7gWJdvMXebxN8h9YSUVZ6vjXo5ZP7MzL6YlaIJj6txOnMvOVcxcnZLwayndknFByWtxfp8p42nBOSCJnnTIfFhhVRhDCI3uLBBjqmAaThJMX71dNkCv3kni%2B3hFSEJnFJdS1RyXwWlaHQtiEKWo8ugd0oAhTN1Qij7mXIcqbaoft%2BqxQYS33yozsS1FqNi2pOd5LQWbg3I83LJlgz5J0cmFf9NDMUdCC2Q5DoWGCrkgbni9doXX3XCfPdtcfsbtvFI9F4OJcRk8JeFlFA%2Bagtk54Tb9MIf0UlFiqRjPt0kNARFwLFZ%2FJL366fRNrIjh%2Bs%2B0E%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2030%20May%202025%2005%3A03%3A20%20GMT&X-Amz-Signature=b1d7b012206965253a810423ee05075844edc5abe980f0818e66b5b82857c544"
$ws.Range("K4").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991153435/1748458785997/c2b7hIofUU-camera-video-segment-driver-1748458790997.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOFTLUGLO%2F20250529%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250529T210320Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEMv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIB2ZTYdObBkz97dPlxhBdWMP%2B%2FWDa8OFSZfyYXQA2ZzgAiA4g56GZxpBRUH4Q4nFWXa6dLWjwViVIhrBPNF6wT0kwirmAwiU%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMuezOu6xnbDig%2BzMOKroD7cPvdDW9pey5Wwt65aZeLUFzA%2FU17M23Vg4GRbV8IjitQm87hihB6MDN%2Fy6kf2QIngVv%2FrPBjQVDos0mRRLFr0oxpIm%2BcLBJKEuUqC9oE%2BnO5oWndifBX1XYlPv688N5yPfu2hpk9D7zxv%2FyEoZLzQP6l3woQXUzYn7xDSlNsHy36DzswM%2FmHAy6tIDwMaI9fz%2BNUs0NdnoYceo8AQIkKwoy4qRmX%2FDVwhVEmazuPIUn4ZP096AlnUDQQr41zrh1bKzzj4qckTCUZTbIOT9VoF%2B02gy5IhJPG8kSvI3Ybrfnz0N%2Ba5Dq6t7INTl1BPr3Ja6Oto31MzTIh9zh6DZvwFIyHbIbSh44iK3zjhYG7JoBQZVHXHdJVEs%2FQmBnONqJoAmGWT8DinWK7v26dUe%2BSUOxDQGRJrGIcOoDlaAlPDHZSF3c2MCtfk4guCBVIdr%2BuJrzSqdG2zwbo4TSopFoXhcLDzf2zNfWJygj9eyvgp9%2FBpBkZhA7gWJdvMXebxN8h9YSUVZ6vjXo5ZP7MzL6YlaIJj6txOnMvOVcxcnZLwayndknFByWtxfp8p42nBOSCJnnTIfFhhVRhDCI3uLBBjqmAaThJMX71dNkCv3kni%2B3hFSEJnFJdS1RyXwWlaHQtiEKWo8ugd0oAhTN1Qij7mXIcqbaoft%2BqxQYS33yozsS1FqNi2pOd5LQWbg3I83LJlgz5J0cmFf9NDMUdCC2Q5DoWGCrkgbni9doXX3XCfPdtcfsbtvFI9F4OJcRk8JeFlFA%2Bagtk54Tb9MIf0UlFiqRjPt0kNARFwLFZ%2FJL366fRNrIjh%2Bs%2B0E%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2030%20May%202025%2005%3A03%3A20%20GMT&X-Amz-Signature=d8b1813029fb0a56b414c6ae0fc6ac31c8c66b6ce91d0517d694530d3a6178db"
$ws.Range("L4").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991153435/1748458785997/eN4S5AQXnX-camera-video-segment-1748458790997.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOFTLUGLO%2F20250529%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250529T210320Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEMv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIB2ZTYdObBkz97dPlxhBdWMP%2B%2FWDa8OFSZfyYXQA2ZzgAiA4g56GZxpBRUH4Q4nFWXa6dLWjwViVIhrBPNF6wT0kwirmAwiU%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMuezOu6xnbDig%2BzMOKroD7cPvdDW9pey5Wwt65aZeLUFzA%2FU17M23Vg4GRbV8IjitQm87hihB6MDN%2Fy6kf2QIngVv%2FrPBjQVDos0mRRLFr0oxpIm%2BcLBJKEuUqC9oE%2BnO5oWndifBX1XYlPv688N5yPfu2hpk9D7zxv%2FyEoZLzQP6l3woQXUzYn7xDSlNsHy36DzswM%2FmHAy6tIDwMaI9fz%2BNUs0NdnoYceo8AQIkKwoy4qRmX%2FDVwhVEmazuPIUn4ZP096AlnUDQQr41zrh1bKzzj4qckTCUZTbIOT9VoF%2B02gy5IhJPG8kSvI3Ybrfnz0N%2Ba5Dq6t7INTl1BPr3Ja6Oto31MzTIh9zh6DZvwFIyHbIbSh44iK3zjhYG7JoBQZVHXHdJVEs%2FQmBnONqJoAmGWT8DinWK7v26dUe%2BSUOxDQGRJrGIcOoDlaAlPDHZSF3c2MCtfk4guCBVIdr%2BuJrzSqdG2zwbo4TSopFoXhcLDzf2zNfWJygj9eyvgp9%2FBpBkZhA7gWJdvMXebxN8h9YSUVZ6vjXo5ZP7MzL6YlaIJj6txOnMvOVcxcnZLwayndknFByWtxfp8p42nBOSCJnnTIfFhhVRhDCI3uLBBjqmAaThJMX71dNkCv3kni%2B3hFSEJnFJdS1RyXwWlaHQtiEKWo8ugd0oAhTN1Qij7mXIcqbaoft%2BqxQYS33yozsS1FqNi2pOd5LQWbg3I83LJlgz5J0cmFf9NDMUdCC2Q5DoWGCrkgbni9doXX3XCfPdtcfsbtvFI9F4OJcRk8JeFlFA%2Bagtk54Tb9MIf0UlFiqRjPt0kNARFwLFZ%2FJL366fRNrIjh%2Bs%2B0E%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2030%20May%202025%2005%3A03%3A20%20GMT&X-Amz-Signature=72fd5bc269ef0f7e2424348b6ac931aa5e39bc40064cdc95935c203dc2b04c5e"
